# Generate Report for Handback
#
# A new handoff/handback round was produced for the source file
# "5826f5c7-56d4-4a44-8440-6cb1cc88e531.md" in both the zh-cn and de-de
# locales. This script updates the "zh-cn" and "de-de" worksheets of the
# localization status report with the new handoff/handback information
# (row 6 of each table), widens the "Error Detail" column so the new,
# longer error text is readable, and adds a hyperlink on the newly
# populated "Latest Target File" cell.

$wb = $excel.ActiveWorkbook

$sourceMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/526582c48a08f35e8397e23aef089b0ec497c0aa/e2e/5826f5c7-56d4-4a44-8440-6cb1cc88e531.md"
$sourceMdDisplay = "5826f5c7-56d4-4a44-8440-6cb1cc88e531.md"

# ---------------------------------------------------------------------
# zh-cn worksheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# New handoff datetime for this round.
$wsZh.Range("H6").Value = "2016-08-26 06:41:42"

# Latest Target File now points at the source markdown file (hyperlink).
$wsZh.Hyperlinks.Add($wsZh.Range("I6"), $sourceMdUrl, "", "", $sourceMdDisplay)

# Latest Handback File uploaded for this round.
$wsZh.Range("J6").Value = "5826f5c7-56d4-4a44-8440-6cb1cc88e531.c75152bbd8fe8d3527a5c005bdc93558bd4beacc.zh-cn.xlf"

# Handback was rejected because the handback file version is stale.
$wsZh.Range("K6").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68917aad07bc037d8e995a203c319cc552a7bb1d/e2e/5826f5c7-56d4-4a44-8440-6cb1cc88e531.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/526582c48a08f35e8397e23aef089b0ec497c0aa/e2e/5826f5c7-56d4-4a44-8440-6cb1cc88e531.md."

# Error detail timestamp.
$wsZh.Range("P6").Value = "2016-08-26 06:41:26"

# Widen the Error Detail column (P) so the long message is readable.
$wsZh.Columns.Item(16).ColumnWidth = 39.14

# ---------------------------------------------------------------------
# de-de worksheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Latest Target File now points at the source markdown file (hyperlink).
$wsDe.Hyperlinks.Add($wsDe.Range("I6"), $sourceMdUrl, "", "", $sourceMdDisplay)

# Latest Handback File uploaded for this round.
$wsDe.Range("J6").Value = "5826f5c7-56d4-4a44-8440-6cb1cc88e531.c75152bbd8fe8d3527a5c005bdc93558bd4beacc.de-de.xlf"

# Handback succeeded for de-de; record the handback datetime.
$wsDe.Range("K6").Value = "2016-08-26 06:41:48"

# Error detail timestamp.
$wsDe.Range("P6").Value = "2016-08-26 06:41:26"

# Widen the Error Detail column (P) so the long message is readable.
$wsDe.Columns.Item(16).ColumnWidth = 39.14
